$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where both Price (D) and Volume(1h) (E) change
$bothChanges = @(
    @{ Row = 2; D = '37.040.62'; E = '  -1.25%  ' }
    @{ Row = 3; D = '2.023.13'; E = '  -2.24%  ' }
    @{ Row = 5; D = '226.19'; E = '  -2.77%  ' }
    @{ Row = 6; D = '0.607'; E = '  -3.74%  ' }
    @{ Row = 8; D = '55.00'; E = '  -4.32%  ' }
    @{ Row = 9; D = '0.380'; E = '  -2.34%  ' }
    @{ Row = 10; D = '0.0788'; E = '  +1.11%  ' }
    @{ Row = 12; D = '2.325.52'; E = '  -2.05%  ' }
    @{ Row = 13; D = '14.31'; E = '  -4.20%  ' }
    @{ Row = 14; D = '20.39'; E = '  -2.68%  ' }
    @{ Row = 15; D = '0.744'; E = '  -2.26%  ' }
    @{ Row = 17; D = '2.032.59'; E = '  -1.28%  ' }
    @{ Row = 18; D = '36.976.06'; E = '  -1.34%  ' }
    @{ Row = 20; D = '68.78'; E = '  -2.51%  ' }
    @{ Row = 21; D = '0.0₃0824'; E = '  -0.66%  ' }
    @{ Row = 22; D = '225.97'; E = '  -0.94%  ' }
    @{ Row = 26; D = '9.24'; E = '  -4.29%  ' }
    @{ Row = 27; D = '165.60'; E = '  -2.31%  ' }
    @{ Row = 29; D = '18.71'; E = '  -3.86%  ' }
    @{ Row = 32; D = '4.47'; E = '  -3.15%  ' }
    @{ Row = 33; D = '0.0617'; E = '  -2.58%  ' }
    @{ Row = 34; D = '4.44'; E = '  -4.21%  ' }
    @{ Row = 35; D = '2.35'; E = '  -5.22%  ' }
    @{ Row = 38; D = '3.16'; E = '  -5.15%  ' }
    @{ Row = 40; D = '17.25'; E = '  +3.21%  ' }
    @{ Row = 43; D = '95.35'; E = '  -4.28%  ' }
    @{ Row = 44; D = '0.0926'; E = '  -3.25%  ' }
    @{ Row = 45; D = '2.79'; E = '  -4.12%  ' }
    @{ Row = 47; D = '7.37'; E = '  +1.61%  ' }
    @{ Row = 50; D = '2.212.77'; E = '  -2.05%  ' }
    @{ Row = 51; D = '3.62'; E = '  -9.30%  ' }
)

foreach ($item in $bothChanges) {
    $dCell = $ws.Cells.Item($item.Row, 4)
    $dCell.Formula = "'" + $item.D
    $dCell.Style = "Normal"
    $eCell = $ws.Cells.Item($item.Row, 5)
    $eCell.Formula = "'" + $item.E
    $eCell.Style = "Normal"
}

# Rows where only Volume(1h) (E) changes
$eOnlyChanges = @(
    @{ Row = 4; E = '  +0.06%  ' }
    @{ Row = 7; E = '  +0.09%  ' }
    @{ Row = 11; E = '  -3.67%  ' }
    @{ Row = 16; E = '  -3.39%  ' }
    @{ Row = 19; E = '  +5.02%  ' }
    @{ Row = 23; E = '  -0.08%  ' }
    @{ Row = 24; E = '  +2.37%  ' }
    @{ Row = 25; E = '  -7.41%  ' }
    @{ Row = 28; E = '  -5.97%  ' }
    @{ Row = 30; E = '  -3.14%  ' }
    @{ Row = 31; E = '  -4.70%  ' }
    @{ Row = 36; E = '  +1.03%  ' }
    @{ Row = 37; E = '  +0.26%  ' }
    @{ Row = 39; E = '  +0.50%  ' }
    @{ Row = 46; E = '  -5.55%  ' }
    @{ Row = 48; E = '  -3.55%  ' }
    @{ Row = 49; E = '  -0.61%  ' }
)

foreach ($item in $eOnlyChanges) {
    $eCell = $ws.Cells.Item($item.Row, 5)
    $eCell.Formula = "'" + $item.E
    $eCell.Style = "Normal"
}

# Rows 41 and 42 swap coin identity (VeChain <-> Maker) with new data
$rowSwaps = @(
    @{ Row = 41; B = 'Maker'; C = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D = '1.488.73'; E = '  +0.58%  ' }
    @{ Row = 42; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.0218'; E = '  -5.44%  ' }
)

foreach ($item in $rowSwaps) {
    $bCell = $ws.Cells.Item($item.Row, 2)
    $bCell.Formula = "'" + $item.B
    $bCell.Style = "Normal"
    $cCell = $ws.Cells.Item($item.Row, 3)
    $cCell.Formula = "'" + $item.C
    $cCell.Style = "Normal"
    $dCell = $ws.Cells.Item($item.Row, 4)
    $dCell.Formula = "'" + $item.D
    $dCell.Style = "Normal"
    $eCell = $ws.Cells.Item($item.Row, 5)
    $eCell.Formula = "'" + $item.E
    $eCell.Style = "Normal"
}

